$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the text of a Range(start,end) with $newText while forcing
# a real rewrite (so embedded <w:proofErr/> markers get dropped) even when
# the visible text doesn't change. We do this by writing a placeholder of
# the same length first, then writing the real text into the (now updated)
# range.
# ---------------------------------------------------------------------------
function Set-RangeTextClean($start, $end, $newText) {
    $r = $d.Range($start, $end)
    $placeholderLen = $r.End - $r.Start
    if ($placeholderLen -gt 0) {
        $r.Text = ("Y" * $placeholderLen)
    }
    $r2 = $d.Range($start, $r.End)
    $r2.Text = $newText
    return $d.Range($start, $r2.End)
}

# ---------------------------------------------------------------------------
# Helper: force a run-split at a given offset pair (so later formatting /
# structural differences keep the text on separate <w:r> runs) without
# changing any visible formatting - toggle Bold on then back off.
# ---------------------------------------------------------------------------
function Split-RunAt($start, $end) {
    $r = $d.Range($start, $end)
    $r.Bold = 1
    $r.Bold = 0
}

# 1) Paragraph 2 : "Não esta salvando rua e numero" merge (do this one first
#    since it sits after paragraph 1 in the body but editing it first keeps
#    paragraph-1 offsets untouched while we still rely on paragraph 2's
#    original offsets).
Set-RangeTextClean 32 85 "Não esta salvando rua e numero e cep, OK TEM Q TESTAR" | Out-Null
Split-RunAt 68 85
Split-RunAt 62 68

# 2) Paragraph 1 : "Cep sem mascara" merge
Set-RangeTextClean 0 31 "Cep sem mascara OK, TEM Q TESTE" | Out-Null
Split-RunAt 15 31

# 3) Paragraph 4 (ATENDIMENTO VETERINARIO ...): append " - Ok" run and move
#    the _GoBack bookmark here (collapsed to the new end of the paragraph).
$p4 = $d.Paragraphs.Item(4).Range
$p4End = $p4.End - 1   # exclude the paragraph mark
$insertion = $d.Range($p4End, $p4End)
$insertion.InsertAfter(" - Ok")
$newEnd = $p4End + 5   # length of " - Ok"
$bm = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null

Write-Host "Done."
